# "Generate Report for Handback"
#
# The workbook tracks localization handoff/handback status for two target
# languages (zh-cn on sheet 2, de-de on sheet 3). This script records that
# both files in each language have now been handed back (in sync with
# en-US), and fills in the "Latest Target File" / "Latest Handback File"
# columns (E/F) that were previously empty, duplicating the same links
# already present in the "Latest Handoff File" column (C) and the
# "Source File Name" column (A). For de-de the handback datetime (column G)
# is also stamped with the actual handback time.

$wb = $excel.ActiveWorkbook

$mdUrl1   = "https://github.com/OpenLocalizationTest/oltest/blob/c162ca4220cfb9c8290489f76192d7792c135d0c/e2e/644510b7-72de-41c5-b953-cfb0bd0c6023.md"
$mdUrl2   = "https://github.com/OpenLocalizationTest/oltest/blob/c162ca4220cfb9c8290489f76192d7792c135d0c/e2e/98644877-9fba-49f4-988d-f79b54b3127b.md"
$mdName1  = "644510b7-72de-41c5-b953-cfb0bd0c6023.md"
$mdName2  = "98644877-9fba-49f4-988d-f79b54b3127b.md"

$statusHandedBack = "Handed back: in sync with en-US"

# ---- Overview sheet (mirrors the per-language Status column) ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value = $statusHandedBack
$wsZh.Range("B3").Value = $statusHandedBack

$zhXlf1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9d9d5df20bf6ccf9e811473a3ce47023f113236f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/644510b7-72de-41c5-b953-cfb0bd0c6023.13c4f007f56ff608b29928e42ebb5e8d2246a2e6.zh-cn.xlf"
$zhXlf2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9d9d5df20bf6ccf9e811473a3ce47023f113236f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/98644877-9fba-49f4-988d-f79b54b3127b.abe0e60c061ec476a0ed333ccd0a434083174562.zh-cn.xlf"
$zhXlfName1 = "644510b7-72de-41c5-b953-cfb0bd0c6023.13c4f007f56ff608b29928e42ebb5e8d2246a2e6.zh-cn.xlf"
$zhXlfName2 = "98644877-9fba-49f4-988d-f79b54b3127b.abe0e60c061ec476a0ed333ccd0a434083174562.zh-cn.xlf"

$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $mdUrl1, "", "", $mdName1)
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhXlf1, "", "", $zhXlfName1)
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), $mdUrl2, "", "", $mdName2)
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhXlf2, "", "", $zhXlfName2)

$wsZh.Range("E2:F3").Font.Underline = 2
$wsZh.Range("E2:F3").Font.Color = 15570276

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = $statusHandedBack
$wsDe.Range("B3").Value = $statusHandedBack

$deXlf1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d7213aeb5ea17aef88c5e81059c3ade7aa8e37d2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/644510b7-72de-41c5-b953-cfb0bd0c6023.13c4f007f56ff608b29928e42ebb5e8d2246a2e6.de-de.xlf"
$deXlf2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d7213aeb5ea17aef88c5e81059c3ade7aa8e37d2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/98644877-9fba-49f4-988d-f79b54b3127b.abe0e60c061ec476a0ed333ccd0a434083174562.de-de.xlf"
$deXlfName1 = "644510b7-72de-41c5-b953-cfb0bd0c6023.13c4f007f56ff608b29928e42ebb5e8d2246a2e6.de-de.xlf"
$deXlfName2 = "98644877-9fba-49f4-988d-f79b54b3127b.abe0e60c061ec476a0ed333ccd0a434083174562.de-de.xlf"

$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $mdUrl1, "", "", $mdName1)
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deXlf1, "", "", $deXlfName1)
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), $mdUrl2, "", "", $mdName2)
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deXlf2, "", "", $deXlfName2)

$wsDe.Range("E2:F3").Font.Underline = 2
$wsDe.Range("E2:F3").Font.Color = 15570276

# de-de handback actually completed -> stamp the real handback datetime
$wsDe.Range("G2").Value = "2016-03-03 13:16:45"
$wsDe.Range("G3").Value = "2016-03-03 13:16:45"
